# Edit script: restructure the stock holdings sheet
#  - Remove the unused "Gain or Loss" / "Percentage" header columns (F, G)
#  - Refresh the figures for the existing holdings (NVAX, SMCI, NTNX)
#  - Append four new holdings (AVXL, VTV, XLK, TSM)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two empty columns (old F: "Gain or Loss", old G: "Percentage")
# that sit between "Current Price" and "Todays Change". Deleting them shifts
# the remaining headers/data left so the layout becomes:
# Stock | Buy Price | Quantity | Invested Amount | Current Price |
# Todays Change | Todays Change Percent | Profit/Loss | Percentage Change | Change Percent
$ws.Range("F1:G1").EntireColumn.Delete()

# Helper data for each holding row: Stock, Buy Price, Quantity, Invested Amount,
# Current Price, Todays Change, Todays Change Percent, Profit/Loss, Percentage Change, Change Percent
$rows = @(
    @("NVAX", 31.73, 844.58, 26798.5234, 8.619999885559082, -0.130000114440918, -1.49, -19518.24389665451, -72.83328116747848, -7.809262812769592),
    @("SMCI", 90.94, 360, 32738.4, 35.06999969482422, 0.6399993896484375, 1.86, -20113.20010986328, -61.43611205759377, -7.809262812769592),
    @("NTNX", 73.16, 175, 12803, 66.69999694824219, -5.650001525878906, -7.81, -1130.500534057617, -8.829965899067536, -7.809262812769592),
    @("AVXL", 6.21, 2065, 12823.65, 9.029999732971191, -0.005000114440917969, -0.06, 5823.299448585511, 45.41062371934286, $null),
    @("VTV", 179.95, 75, 13496.25, 181.2599945068359, -0.0800018310546875, -0.04, 98.24958801269617, 0.7279769418371486, $null),
    @("XLK", 234.42, 135, 31646.7, 231.5800018310547, -3.17999267578125, -1.35, -383.3997528076155, -1.211499944094062, $null),
    @("TSM", 193.64, 12, 2323.68, 181.1900024414062, -2.649993896484375, -1.44, -149.3999707031248, -6.429455463020934, $null)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    if ($null -ne $row[9]) {
        $ws.Cells.Item($r, 10).Value = $row[9]
    }
    $r = $r + 1
}
